# Auto-generated Excel COM-interop script
# Applies numeric data updates to the Anima Profits leve-crafting sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2314.2856
$ws.Range("I51").Value = 1266.6666
$ws.Range("J51").Value = 3100
$ws.Range("K51").Value = 1266.6666
$ws.Range("L51").Value = 3100
$ws.Range("M51").Value = -782.6666
$ws.Range("N51").Value = -4068
$ws.Range("H88").Value = 61380.8
$ws.Range("I88").Value = 1750
$ws.Range("J88").Value = 101134.664
$ws.Range("K88").Value = 1750
$ws.Range("L88").Value = 101134.664
$ws.Range("M88").Value = -1344
$ws.Range("N88").Value = -101946.664
$ws.Range("H91").Value = 61380.8
$ws.Range("I91").Value = 1750
$ws.Range("J91").Value = 101134.664
$ws.Range("K91").Value = 1750
$ws.Range("L91").Value = 101134.664
$ws.Range("M91").Value = -346
$ws.Range("N91").Value = -103942.664
$ws.Range("H116").Value = 6132.04
$ws.Range("I116").Value = 10512.917
$ws.Range("J116").Value = 2088.1538
$ws.Range("K116").Value = 10512.917
$ws.Range("L116").Value = 2088.1538
$ws.Range("M116").Value = -7070.916999999999
$ws.Range("N116").Value = -8972.1538
$ws.Range("H132").Value = 1638.225
$ws.Range("I132").Value = 1444.0405
$ws.Range("J132").Value = 4033.1667
$ws.Range("K132").Value = 4332.1215
$ws.Range("L132").Value = 12099.5001
$ws.Range("M132").Value = -1802.1215
$ws.Range("N132").Value = -17159.5001
$ws.Range("H138").Value = 1215.86
$ws.Range("I138").Value = 554.21277
$ws.Range("J138").Value = 1802.6038
$ws.Range("K138").Value = 1662.63831
$ws.Range("L138").Value = 5407.811400000001
$ws.Range("M138").Value = 3477.36169
$ws.Range("N138").Value = -15687.8114
$ws.Range("H141").Value = 2186.9285
$ws.Range("I141").Value = 870.55554
$ws.Range("J141").Value = 7572.091
$ws.Range("K141").Value = 2611.66662
$ws.Range("L141").Value = 22716.273
$ws.Range("M141").Value = 2568.33338
$ws.Range("N141").Value = -33076.273
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 725650.4
$ws.Range("I32").Value = 850410.6
$ws.Range("K32").Value = 850410.6
$ws.Range("M32").Value = -850123.6
$ws.Range("H132").Value = 2559.1187
$ws.Range("I132").Value = 2608.5
$ws.Range("J132").Value = 2491.96
$ws.Range("K132").Value = 7825.5
$ws.Range("L132").Value = 7475.88
$ws.Range("M132").Value = -5295.5
$ws.Range("N132").Value = -12535.88
$ws.Range("H137").Value = 19857.143
$ws.Range("J137").Value = 18000
$ws.Range("L137").Value = 18000
$ws.Range("N137").Value = -28200
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1279.3077
$ws.Range("I99").Value = 977.25
$ws.Range("J99").Value = 1762.6
$ws.Range("K99").Value = 977.25
$ws.Range("L99").Value = 1762.6
$ws.Range("M99").Value = 520.75
$ws.Range("N99").Value = -4758.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4077.1594
$ws.Range("I31").Value = 1235
$ws.Range("J31").Value = 7561.0967
$ws.Range("K31").Value = 1235
$ws.Range("L31").Value = 7561.0967
$ws.Range("M31").Value = -940
$ws.Range("N31").Value = -8151.0967
$ws.Range("H34").Value = 4077.1594
$ws.Range("I34").Value = 1235
$ws.Range("J34").Value = 7561.0967
$ws.Range("K34").Value = 1235
$ws.Range("L34").Value = 7561.0967
$ws.Range("M34").Value = -1033
$ws.Range("N34").Value = -7965.0967
$ws.Range("H58").Value = 1066.2554
$ws.Range("I58").Value = 786.64514
$ws.Range("J58").Value = 1608
$ws.Range("K58").Value = 786.64514
$ws.Range("L58").Value = 1608
$ws.Range("M58").Value = -583.64514
$ws.Range("N58").Value = -2014
$ws.Range("H132").Value = 3087834.2
$ws.Range("I132").Value = 1234.3513
$ws.Range("K132").Value = 3703.0539
$ws.Range("M132").Value = -1173.0539
$ws.Range("H136").Value = 1066.2554
$ws.Range("I136").Value = 786.64514
$ws.Range("J136").Value = 1608
$ws.Range("K136").Value = 2359.93542
$ws.Range("L136").Value = 4824
$ws.Range("M136").Value = 190.0645800000002
$ws.Range("N136").Value = -9924
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1161.4147
$ws.Range("J5").Value = 2311.7646
$ws.Range("L5").Value = 6935.293799999999
$ws.Range("N5").Value = -7159.293799999999
$ws.Range("H122").Value = 2776.8262
$ws.Range("J122").Value = 6485.6113
$ws.Range("L122").Value = 58370.50169999999
$ws.Range("N122").Value = -63270.50169999999
$ws.Range("H135").Value = 1161.4147
$ws.Range("J135").Value = 2311.7646
$ws.Range("L135").Value = 20805.8814
$ws.Range("N135").Value = -25875.8814
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 251
$ws.Range("J93").Value = 251
$ws.Range("L93").Value = 251
$ws.Range("N93").Value = -3995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 24135.666
$ws.Range("J24").Value = 24135.666
$ws.Range("L24").Value = 24135.666
$ws.Range("N24").Value = -24821.666
$ws.Range("H132").Value = 2364.7273
$ws.Range("I132").Value = 2019.1072
$ws.Range("J132").Value = 4300.2
$ws.Range("K132").Value = 6057.321599999999
$ws.Range("L132").Value = 12900.6
$ws.Range("M132").Value = -3527.321599999999
$ws.Range("N132").Value = -17960.6
$ws.Range("H136").Value = 4168329
$ws.Range("I136").Value = 1454.7693
$ws.Range("J136").Value = 11906809
$ws.Range("K136").Value = 4364.3079
$ws.Range("L136").Value = 35720427
$ws.Range("M136").Value = -1814.3079
$ws.Range("N136").Value = -35725527
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 2693333.2
$ws.Range("H122").Value = 2774.0889
$ws.Range("I122").Value = 2600.9355
$ws.Range("J122").Value = 3157.5
$ws.Range("K122").Value = 7802.806500000001
$ws.Range("L122").Value = 9472.5
$ws.Range("M122").Value = -5352.806500000001
$ws.Range("N122").Value = -14372.5
$ws.Range("H132").Value = 4488674
$ws.Range("I132").Value = 1496.659
$ws.Range("J132").Value = 13890378
$ws.Range("K132").Value = 4489.977000000001
$ws.Range("L132").Value = 41671134
$ws.Range("M132").Value = -1959.977000000001
$ws.Range("N132").Value = -41676194
$ws.Range("H136").Value = 1826.0505
$ws.Range("I136").Value = 1781.2924
$ws.Range("K136").Value = 5343.8772
$ws.Range("M136").Value = -2793.8772
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 17634.166
$ws.Range("I93").Value = 17634.166
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 17634.166
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -16386.166
$ws.Range("N93").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

Write-Host "Applied all leve profit updates"
